# Commit: "changing all references from TEO to PEO"
#
# 1) The datetimeFigureOut date field ("7/13/2013" -> "7/14/2013") that is
#    rendered on every slide actually lives on the slide layouts + slide
#    master (this deck only has one slide, but eleven layouts + the master
#    all carry the placeholder), so update it there.
# 2) On the single content slide, three titles get split into multiple runs
#    (mixed-case retyping of "Array"/"Atomic Segment"/"Order" into lower
#    case second words) and two digit textboxes lose a redundant trailing
#    endParaRPr.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($container)
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {}
        if ($isDatePh) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "7/13/2013") {
                $tr.Text = "7/14/2013"
            }
        }
    }
}

# --- 1) Date placeholder on every slide layout + the slide master ---
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j)
}
Update-DatePlaceholder $p.SlideMaster

# --- helper to find a shape by name on the content slide ---
function Get-ShapeByName {
    param($slide, [string]$name)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

$s = $p.Slides.Item(1)

# --- 2) "Buffer Array" -> "Buffer " + "array" ---
$shp7 = Get-ShapeByName $s "TextBox 7"
$tr7 = $shp7.TextFrame.TextRange
$tr7.Characters(8, 5).Text = "array"

# --- 3) Drop the redundant trailing endParaRPr on the "2" / "3" boxes ---
$shp10 = Get-ShapeByName $s "TextBox 10"
$tr10 = $shp10.TextFrame.TextRange
$tr10.Delete()
$tr10.InsertBefore("2")

$shp11 = Get-ShapeByName $s "TextBox 11"
$tr11 = $shp11.TextFrame.TextRange
$tr11.Delete()
$tr11.InsertBefore("3")

# --- 4) "Buffer Atomic Segment" -> "Buffer " + "atomic " + "s" + "egment" ---
$shp15 = Get-ShapeByName $s "TextBox 15"
$tr15 = $shp15.TextFrame.TextRange
$tr15.Characters(8, 7).Text = "atomic "
$tr15.Characters(15, 1).Text = "s"

# --- 5) "Persist Order" -> "Persist " + "order" ---
$shp19 = Get-ShapeByName $s "TextBox 19"
$tr19 = $shp19.TextFrame.TextRange
$tr19.Characters(9, 5).Text = "order"
